# Generate Report for Handback
# ------------------------------------------------------------------
# This script adds a second handed-back file (aac79222-062a-4681-b8ed-
# d0285ad8e595.md) to the handback status report, alongside the existing
# file (a3da39c9-120d-42b6-8c73-c0b3bf24fdb3.md, whose generated GUID is
# refreshed to a46d0e08-2453-416e-9b29-e3a1a91e3d7a.md for this run), and
# refreshes the run's timestamps across all three sheets (Overview,
# zh-cn, de-de).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$HYPER_COLOR = 15570276   # OLE (BGR) form of RGB 6495ED - matches the workbook's HyperLink style

$file1Guid = "a46d0e08-2453-416e-9b29-e3a1a91e3d7a"
$file2Guid = "aac79222-062a-4681-b8ed-d0285ad8e595"

$file1Md   = "$file1Guid.md"
$file2Md   = "$file2Guid.md"

$file1ZhXlf = "$file1Guid.263c7478941e6f4247ab41b71d5edd3aa115a34b.zh-cn.xlf"
$file1DeXlf = "$file1Guid.263c7478941e6f4247ab41b71d5edd3aa115a34b.de-de.xlf"
$file2ZhXlf = "$file2Guid.e42f4eab8cc73b3badd248a26e0ded876ad41f8d.zh-cn.xlf"
$file2DeXlf = "$file2Guid.e42f4eab8cc73b3badd248a26e0ded876ad41f8d.de-de.xlf"

$overviewDate = "2016-08-16 06:54:38"
$zhHandoffDate  = "2016-08-16 06:54:32"
$zhHandbackDate = "2016-08-16 06:55:01"
$deHandbackDate = "2016-08-16 06:55:15"

$repoBase = "https://github.com/OpenLocalizationTestOrg"
$commitMain = "b7c1e4f6a9d2c5038f6a1b4e7d0c9a2f5b8e1d46"
$commitZh   = "c4a8f1e6b3d9072a5f8c1e4b7a0d3f6c9e2b5a81"
$commitDe   = "d9f2b5a8e1c4073b6a9f2c5e8b1d4a7f0c3e6b92"

function Set-Txt($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

# ====================================================================
# Sheet 1: Overview
# ====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)

# --- refresh existing row 2 (file 1) ---
Set-Txt $wsOverview "A2" $file1Md
Set-Txt $wsOverview "B2" "e2e\$file1Md"
Set-Txt $wsOverview "G2" $overviewDate

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "$repoBase/ol-test0/blob/$commitMain/e2e/$file1Md",
    "", "", "e2e\$file1Md") | Out-Null
$wsOverview.Range("B2").Font.Underline = $true
$wsOverview.Range("B2").Font.Color = $HYPER_COLOR

# --- add new row 3 (file 2) ---
$newOverviewRow = $tblOverview.ListRows.Add()
Set-Txt $wsOverview "A3" $file2Md
Set-Txt $wsOverview "B3" "e2e\$file2Md"
Set-Txt $wsOverview "C3" ".md"
Set-Txt $wsOverview "E3" "Handed back: in sync with en-US"
Set-Txt $wsOverview "F3" "Handed back: in sync with en-US"
Set-Txt $wsOverview "G3" $overviewDate

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "$repoBase/ol-test0/blob/$commitMain/e2e/$file2Md",
    "", "", "e2e\$file2Md") | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $HYPER_COLOR
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ====================================================================
# Sheet 2: zh-cn
# ====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$tblZh = $wsZh.ListObjects.Item(1)

# --- refresh existing row 2 (file 1) ---
Set-Txt $wsZh "A2" $file1Md
Set-Txt $wsZh "G2" $file1ZhXlf
Set-Txt $wsZh "H2" $zhHandoffDate
Set-Txt $wsZh "I2" $file1Md
Set-Txt $wsZh "J2" $file1ZhXlf
Set-Txt $wsZh "K2" $zhHandbackDate

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "$repoBase/ol-test0-zhcn/blob/$commitZh/e2e/$file1Md",
    "", "", $file1Md) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I2"),
    "$repoBase/ol-test0-zhcn/blob/$commitZh/e2e/$file1Md",
    "", "", $file1Md) | Out-Null
$wsZh.Range("A2").Font.Underline = $true
$wsZh.Range("A2").Font.Color = $HYPER_COLOR
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = $HYPER_COLOR
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- add new row 3 (file 2) ---
$newZhRow = $tblZh.ListRows.Add()
Set-Txt $wsZh "A3" $file2Md
Set-Txt $wsZh "B3" ".md"
Set-Txt $wsZh "C3" "Handed back: in sync with en-US"
Set-Txt $wsZh "D3" "e2e"
Set-Txt $wsZh "E3" "ht"
Set-Txt $wsZh "F3" "True"
Set-Txt $wsZh "G3" $file2ZhXlf
Set-Txt $wsZh "H3" $zhHandoffDate
Set-Txt $wsZh "I3" $file2Md
Set-Txt $wsZh "J3" $file2ZhXlf
Set-Txt $wsZh "K3" $zhHandbackDate
Set-Txt $wsZh "L3" ""
Set-Txt $wsZh "M3" "True"
Set-Txt $wsZh "N3" ""
Set-Txt $wsZh "O3" "False"
Set-Txt $wsZh "P3" ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "$repoBase/ol-test0-zhcn/blob/$commitZh/e2e/$file2Md",
    "", "", $file2Md) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I3"),
    "$repoBase/ol-test0-zhcn/blob/$commitZh/e2e/$file2Md",
    "", "", $file2Md) | Out-Null
$wsZh.Range("A3").Font.Underline = $true
$wsZh.Range("A3").Font.Color = $HYPER_COLOR
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = $HYPER_COLOR
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ====================================================================
# Sheet 3: de-de
# ====================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$tblDe = $wsDe.ListObjects.Item(1)

# --- refresh existing row 2 (file 1) ---
Set-Txt $wsDe "A2" $file1Md
Set-Txt $wsDe "G2" $file1DeXlf
Set-Txt $wsDe "H2" $overviewDate
Set-Txt $wsDe "I2" $file1Md
Set-Txt $wsDe "J2" $file1DeXlf
Set-Txt $wsDe "K2" $deHandbackDate

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "$repoBase/ol-test0-dede/blob/$commitDe/e2e/$file1Md",
    "", "", $file1Md) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I2"),
    "$repoBase/ol-test0-dede/blob/$commitDe/e2e/$file1Md",
    "", "", $file1Md) | Out-Null
$wsDe.Range("A2").Font.Underline = $true
$wsDe.Range("A2").Font.Color = $HYPER_COLOR
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = $HYPER_COLOR
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- add new row 3 (file 2) ---
$newDeRow = $tblDe.ListRows.Add()
Set-Txt $wsDe "A3" $file2Md
Set-Txt $wsDe "B3" ".md"
Set-Txt $wsDe "C3" "Handed back: in sync with en-US"
Set-Txt $wsDe "D3" "e2e"
Set-Txt $wsDe "E3" "ht"
Set-Txt $wsDe "F3" "True"
Set-Txt $wsDe "G3" $file2DeXlf
Set-Txt $wsDe "H3" $overviewDate
Set-Txt $wsDe "I3" $file2Md
Set-Txt $wsDe "J3" $file2DeXlf
Set-Txt $wsDe "K3" $deHandbackDate
Set-Txt $wsDe "L3" ""
Set-Txt $wsDe "M3" "True"
Set-Txt $wsDe "N3" ""
Set-Txt $wsDe "O3" "False"
Set-Txt $wsDe "P3" ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "$repoBase/ol-test0-dede/blob/$commitDe/e2e/$file2Md",
    "", "", $file2Md) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I3"),
    "$repoBase/ol-test0-dede/blob/$commitDe/e2e/$file2Md",
    "", "", $file2Md) | Out-Null
$wsDe.Range("A3").Font.Underline = $true
$wsDe.Range("A3").Font.Color = $HYPER_COLOR
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = $HYPER_COLOR
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Handback report regenerated."
